$d = $word.ActiveDocument

# Locate the paragraph containing the "Please read Using Pressure Canners..."
# sentence so we can remove it along with the blank paragraph that follows it.
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute( `
    "Please read Using Pressure Canners before beginning. If this is your first time canning, it is recommended that you read Principles of Home Canning.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Paragraph 1 is the sentence paragraph itself; paragraph 2 is the
    # following blank ("\n") paragraph. Removing the range spanning both
    # paragraph marks deletes both paragraphs entirely.
    $startPos = $findRange.Paragraphs(1).Range.Start
    $endPos = $findRange.Paragraphs(2).Range.End
    $d.Range($startPos, $endPos).Delete()
}
